$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update November 2025 (row 24) stats with the newly reported figures.
$ws.Range("B24").Value = 6339
$ws.Range("C24").Value = 999
$ws.Range("D24").Value = 5940212
$ws.Range("E24").Value = 937.0897617920808
$ws.Range("F24").Value = 8.063416297306514
$ws.Range("G24").Value = 3.523316062176174
$ws.Range("H24").Value = 25.83669945082927
